$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the species data between row 2 and row 4 (columns A, B, E, F, G, H),
# and update B3's value, per the target edit.

# Save row 2 values before overwriting
$row2A = $ws.Range("A2").Value2
$row2B = $ws.Range("B2").Value2
$row2E = $ws.Range("E2").Value2
$row2F = $ws.Range("F2").Value2
$row2G = $ws.Range("G2").Value2
$row2H = $ws.Range("H2").Value2

# Save row 4 values
$row4A = $ws.Range("A4").Value2
$row4B = $ws.Range("B4").Value2
$row4E = $ws.Range("E4").Value2
$row4F = $ws.Range("F4").Value2
$row4G = $ws.Range("G4").Value2
$row4H = $ws.Range("H4").Value2

# Write row 4's original values into row 2
$ws.Range("A2").Value = $row4A
$ws.Range("B2").Value = $row4B
$ws.Range("E2").Value = $row4E
$ws.Range("F2").Value = $row4F
$ws.Range("G2").Value = $row4G
$ws.Range("H2").Value = $row4H

# Write row 2's original values into row 4
$ws.Range("A4").Value = $row2A
$ws.Range("B4").Value = $row2B
$ws.Range("E4").Value = $row2E
$ws.Range("F4").Value = $row2F
$ws.Range("G4").Value = $row2G
$ws.Range("H4").Value = $row2H

# Update B3 value
$ws.Range("B3").Value = 89769
